$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.171.52"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "'2.049.42"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D5").Value = "'251.04"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'64.39"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.390"
$ws.Range("E9").Value = "  +5.36%  "
$ws.Range("D10").Value = "'58.53"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  +7.34%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("E14").Value = "  +16.97%  "
$ws.Range("D15").Value = "'14.53"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "'2.350.94"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "'5.59"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").Value = "'2.053.24"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "'37.111.39"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").Value = "'72.84"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "'0.0₃0891"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").Value = "'237.65"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  -5.61%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "'160.05"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").Value = "'20.32"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("D30").Value = "'0.135"
$ws.Range("E30").Value = "  +23.95%  "
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'0.0625"
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("D35").Value = "'4.57"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'2.39"
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("D37").Value = "'6.35"
$ws.Range("E37").Value = "  +10.21%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").Value = "'3.03"
$ws.Range("E40").Value = "  +26.64%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.26"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.100"
$ws.Range("E42").Value = "  -8.06%  "
$ws.Range("D43").Value = "'2.99"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'17.27"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'94.25"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'7.79"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "'1.374.62"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'2.237.54"
$ws.Range("E51").Value = "  +1.67%  "
